$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: text-only swaps (Find/Replace scoped to a single paragraph's Range
# so there is zero risk of cross-matching other paragraphs). None of these
# change paragraph count, so paragraph indices stay stable throughout.
# ---------------------------------------------------------------------------

# Docente(s) paragraph (index 9): two runs swap to the old "Objetivos" body
# and the old "Programa" body.
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Find.Execute(
    "5840897 - Clodoaldo Saron", $true, $false, $false, $false, $false, $true, 1, $false,
    "Abordar conceitos fundamentais sobre materiais poliméricos, envolvendo o histórico de desenvolvimento, organização do setor produtivo, formas de obtenção, estrutura química e física destes materiais e respectiva relação com propriedades térmicas, mecânicas, reológicas, ópticas, elétricas, etc. ^lPermitir que o aluno tenha uma visão clara sobre estrutura, propriedades e aplicações de polímeros termoplásticos, termorrígidos e elastômeros, bem como as propriedades destes materiais podem ser modificadas com o uso de aditivos.^lCapacitar o aluno com conhecimentos para que possa cursar outras disciplinas na área de materiais poliméricos.",
    2) | Out-Null

$p9.Range.Find.Execute(
    "1033242 - Fábio Herbst Florenzano", $true, $false, $false, $false, $false, $true, 1, $false,
    "Introdução: Desenvolvimento dos materiais poliméricos, organização da cadeia produtiva, formas de obtenção, nomenclatura, arquitetura molecular e estrutura configuracional. Estado sólido: amorfo, cristalino e elastomérico. Estrutura e propriedades. Thermoplásticos: estrutura, propriedades e aplicações. Elastômeros: estrutura, propriedades e aplicações. Resinas termorrígidas: estrutura, propriedades e aplicações. Propriedades mecânicas dos polímeros: comportamento à tração, impacto, flexão e fluência. Aditivos para polímeros: classes e aplicações. Viagem Didática complementar",
    2) | Out-Null

# Programa paragraph (index 14): becomes the old "Critério" evaluation text.
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Find.Execute(
    "Introdução: Desenvolvimento dos materiais poliméricos, organização da cadeia produtiva, formas de obtenção, nomenclatura, arquitetura molecular e estrutura configuracional. Estado sólido: amorfo, cristalino e elastomérico. Estrutura e propriedades. Thermoplásticos: estrutura, propriedades e aplicações. Elastômeros: estrutura, propriedades e aplicações. Resinas termorrígidas: estrutura, propriedades e aplicações. Propriedades mecânicas dos polímeros: comportamento à tração, impacto, flexão e fluência. Aditivos para polímeros: classes e aplicações. Viagem Didática complementar",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação.",
    2) | Out-Null

# Avaliação paragraph (index 16): three plain runs (after the bold labels)
# swap to the old "Norma de recuperação" text, the Bibliografia content, and
# the first docente's name.
$p16 = $d.Paragraphs.Item(16)

$p16.Range.Find.Execute(
    "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5.",
    2) | Out-Null

$p16.Range.Find.Execute(
    "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SIMAL, A. L. Estrutura e Propriedades dos Polímeros, EduFSCar, São Carlos, 2002.^lSPERLING, L. H. Introduction to Physical Polymer Science, New York, John Wiley & Sans, 1985.^lBRYDSON, J. A. Rubbery Materials and Their Compounds, Elsevier, London, 1988.^lRabello, M. S. Aditivação de Polímeros, Artiliber, São Paulo, 2004.^lHARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill Inc, 1992 S. V. CANEVAROLO Jr. Técnicas de Caracterização de Polímeros. São Paulo: Editora Artliber, 2005. MANRICH, S. Processamento de Termoplásticos. Editora Artliber, 2005. NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997. MANO, E. B.; MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000. TURI, E. A. Thermal Characterization of Polymeric Materials. New York: Academic Press, 1981. NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997.MANO, E. B.; MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000",
    2) | Out-Null

$p16.Range.Find.Execute(
    "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5840897 - Clodoaldo Saron",
    2) | Out-Null

# Bibliografia paragraph (index 18): whole multi-run body collapses into the
# second docente's name (a single run, matching the target structure).
$p18 = $d.Paragraphs.Item(18)
$p18.Range.Text = "1033242 - Fábio Herbst Florenzano"

# ---------------------------------------------------------------------------
# Step 2: paragraph merges (each reduces the total paragraph count by one).
# Handle the higher-index merge first so the lower paragraph's index (6/7)
# is unaffected by the shift caused by the first merge.
# ---------------------------------------------------------------------------

# Programa resumido (indices 11 and 12): merge away the empty/italic
# paragraph break and collapse into a single run with the Avaliação "Método"
# text.
$p11 = $d.Paragraphs.Item(11)
$mergeMark1 = $d.Range($p11.Range.End - 1, $p11.Range.End)
$mergeMark1.Delete()
$d.Paragraphs.Item(11).Range.Text = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula"

# Objetivos (indices 6 and 7): merge away the trailing empty/italic
# paragraph and collapse into a single run with the old "Programa resumido"
# text.
$p6 = $d.Paragraphs.Item(6)
$mergeMark2 = $d.Range($p6.Range.End - 1, $p6.Range.End)
$mergeMark2.Delete()
$d.Paragraphs.Item(6).Range.Text = "Histórico dos polímeros, Estrutura e Propriedades de Materiais Poliméricos, Polímeros Termoplásticos, Elastômeros, Polímeros Termorrígidos e Aditivos para Polímeros"
